$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.946.92'
$ws.Range('E2').Value = '  -1.84%  '
$ws.Range('D3').Value = '3.378.64'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''573.77'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.55%  '
$ws.Range('D6').Value = '''136.36'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.75%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.377.34'
$ws.Range('E8').Value = '  -1.01%  '
$ws.Range('E9').Value = '  -1.51%  '
$ws.Range('E10').Value = '  +2.01%  '
$ws.Range('D11').Value = '''0.122'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -4.12%  '
$ws.Range('E12').Value = '  -2.92%  '
$ws.Range('D13').Value = '3.954.70'
$ws.Range('E13').Value = '  -1.04%  '
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '3.384.59'
$ws.Range('E15').Value = '  -0.85%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '''0.0000172'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -4.17%  '
$ws.Range('D17').Value = '''25.41'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('D18').Value = '61.134.41'
$ws.Range('E18').Value = '  -1.75%  '
$ws.Range('D19').Value = '''13.82'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.78%  '
$ws.Range('D20').Value = '''5.74'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.77%  '
$ws.Range('D21').Value = '''9.33'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.61%  '
$ws.Range('D22').Value = '''376.50'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -4.92%  '
$ws.Range('D23').Value = '3.518.46'
$ws.Range('E23').Value = '  -0.95%  '
$ws.Range('E24').Value = '  -3.13%  '
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('E26').Value = '  -4.71%  '
$ws.Range('D27').Value = '''70.94'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.07%  '
$ws.Range('D28').Value = '''0.181'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +11.92%  '
$ws.Range('D29').Value = '''1.63'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.88%  '
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('D31').Value = '''7.39'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -4.49%  '
$ws.Range('D32').Value = '''8.06'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.63%  '
$ws.Range('D33').Value = '''2.14'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.30%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').Value = '''23.52'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.39%  '
$ws.Range('E36').Value = '  -4.91%  '
$ws.Range('E37').Value = '  -3.58%  '
$ws.Range('D38').Value = '''6.78'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.59%  '
$ws.Range('D39').Value = '''164.30'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.35%  '
$ws.Range('D40').Value = '''0.0754'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -4.89%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = '''1.00'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '''25.35'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.23%  '
$ws.Range('D43').Value = '''0.773'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.08%  '
$ws.Range('D44').Value = '''41.70'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.28%  '
$ws.Range('D45').Value = '''1.69'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -6.20%  '
$ws.Range('D46').Value = '''1.19'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -6.60%  '
$ws.Range('D47').Value = '''4.34'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.89%  '
$ws.Range('D48').Value = '2.481.73'
$ws.Range('E48').Value = '  +4.99%  '
$ws.Range('D49').Value = '''6.78'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.29%  '
$ws.Range('D50').Value = '''22.81'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.84%  '
$ws.Range('E51').Value = '  +2.94%  '
